$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that follows the H1 title
#    paragraph at the top of the document.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$null = $metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new paragraph ("Play Fu Dao Le for Free: Exciting Slot
#    Game Review" in bold) right before the final "Prompt: ..." / image
#    description paragraph at the very end of the document.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($count - 1)
$null = $prevPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fu Dao Le for Free: Exciting Slot Game Review</w:t></w:r></w:p>'
$null = $newPara.Range.InsertXML($xml)

# ------------------------------------------------------------------
# 3) Replace the text of the last paragraph (formerly the image-prompt
#    text) with the meta-description text, keeping its formatting.
# ------------------------------------------------------------------
$oldText = "Prompt: Create a cartoon-style feature image for Fu Dao Le that features a happy Maya Warrior wearing glasses. The image should have a vibrant Chinese theme with traditional elements such as red and gold in the background. The warrior should be depicted as victorious and holding a pile of gold coins with a big smile on their face. It should also include the game title " + [char]34 + "Fu Dao Le" + [char]34 + " in bold letters at the top of the image."
$newText = "Looking for an engaging and authentic Chinese-themed slot game? Read our review of Fu Dao Le and play for free!"

$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
